$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.756.56'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.855.15'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '''313.93'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '''0.4331'
$ws.Range("E7").Value = '  +1.65%  '
$ws.Range("D8").Value = '''0.3666'
$ws.Range("E8").Value = '  -0.64%  '
$ws.Range("D9").Value = '''45.08'
$ws.Range("E9").Value = '  +1.56%  '
$ws.Range("D10").Value = '''0.07351'
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").Value = '''0.8828'
$ws.Range("E11").Value = '  -2.57%  '
$ws.Range("D12").Value = '''20.85'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '1.872.82'
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("D14").Value = '''5.368'
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").Value = '''6.554'
$ws.Range("E15").Value = '  -1.08%  '
$ws.Range("D16").Value = '''0.06946'
$ws.Range("E16").Value = '  +1.55%  '
$ws.Range("D17").Value = '''1.004'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("E18").Value = '  +3.21%  '
$ws.Range("D19").Value = '''0.000009096'
$ws.Range("E19").Value = '  +2.38%  '
$ws.Range("D20").Value = '''1.004'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = '''15.44'
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").Value = '28.002.12'
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("D23").Value = '''4.997'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = '''10.49'
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").Value = '2.163.12'
$ws.Range("E25").Value = '  +3.96%  '
$ws.Range("D26").Value = '''1.994'
$ws.Range("E26").Value = '  -2.55%  '
$ws.Range("D27").Value = '''156.14'
$ws.Range("E27").Value = '  +1.31%  '
$ws.Range("D28").Value = '''18.71'
$ws.Range("E28").Value = '  +2.36%  '
$ws.Range("D29").Value = '''5.348'
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("D30").Value = '''121.19'
$ws.Range("E30").Value = '  +8.54%  '
$ws.Range("D31").Value = '''1.871'
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("D32").Value = '''0.08957'
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("D33").Value = '''0.7702'
$ws.Range("E33").Value = '  -2.28%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''4.579'
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").Value = '''2.979'
$ws.Range("E35").Value = '  +2.41%  '
$ws.Range("D36").Value = '''1.127'
$ws.Range("E36").Value = '  +2.99%  '
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '''1.001'
$ws.Range("E37").Value = '  +0.45%  '
$ws.Range("D38").Value = '''1.113'
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.05454'
$ws.Range("E39").Value = '  +0.89%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '''0.01949'
$ws.Range("E40").Value = '  +0.89%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '''2.863'
$ws.Range("E41").Value = '  -3.95%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.5122'
$ws.Range("E42").Value = '  +0.71%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.1666'
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '''6.756'
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '''8.418'
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("D46").Value = '''10.45'
$ws.Range("E46").Value = '  +0.19%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.06551'
$ws.Range("E47").Value = '  -1.71%  '
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '''0.4694'
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = '''104.99'
$ws.Range("E49").Value = '  -1.01%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = '''1.001'
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '''1.632'
$ws.Range("E51").Value = '  -0.84%  '
